$d = $word.ActiveDocument
Write-Output "count=$($d.Styles.Count)"
for ($i = 1; $i -le $d.Styles.Count; $i++) {
  $s = $d.Styles.Item($i)
  Write-Output "$i : $($s.NameLocal) type=$($s.Type)"
}
